# Edit the SmartArt ("gameFlowChart") diagram on slide 1.
#
# The node that currently reads "Runs momo gen twice. " is being split into a
# new, separate death-node class (GameNodeWithGameDeathNode) rather than being
# reused as part of GameNodeWithMoMo, so its trailing " twice." wording is
# trimmed down to just "Runs momo gen".
#
# This text lives inside the diagram's data model (dgm:t) and is mirrored
# automatically into the cached drawing (dsp:txBody) by the host when the
# SmartArt node's text is updated through the object model.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)

$oldText = "Runs momo gen twice. "
$newText = "Runs momo gen"

for ($si = 1; $si -le $s.Shapes.Count; $si++) {
    $shp = $s.Shapes.Item($si)

    if (-not $shp.HasSmartArt) {
        continue
    }

    $smartArt = $shp.SmartArt
    $nodes = $smartArt.AllNodes

    for ($ni = 1; $ni -le $nodes.Count; $ni++) {
        $node = $nodes.Item($ni)
        $tr = $node.TextFrame2.TextRange

        if ($tr.Text -eq $oldText) {
            $tr.Text = $newText
        }
    }
}
